# Updated symbol list on Tue Jan  3 07:50:59 UTC 2023 with GitHub Actions
#
# This updates the cryptocurrency price/volume table on Sheet1.
# All Price/Volume(1h) cells in this sheet are stored as plain TEXT
# (e.g. "246.18", "0.26%"), not numbers, so a leading apostrophe is
# used when assigning values that look numeric to force Excel to keep
# them as text instead of auto-converting to Number/Percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (BNB) ---
$ws.Range("D2").Value = "'246.25"
$ws.Range("E2").Value = "'0.30%"

# --- Row 3 (OKB) ---
$ws.Range("D3").Value = "'29.95"
$ws.Range("E3").Value = "'0.16%"

# --- Row 4 (HuobiToken) ---
$ws.Range("D4").Value = "'5.164"
$ws.Range("E4").Value = "'0.52%"

# --- Row 5 (Cronos) ---
$ws.Range("D5").Value = "'0.05798"
$ws.Range("E5").Value = "'1.15%"

# --- Row 6 (KuCoinToken) ---
$ws.Range("D6").Value = "'6.678"
$ws.Range("E6").Value = "'1.57%"

# --- Row 7 (GateToken) ---
$ws.Range("D7").Value = "'3.215"
$ws.Range("E7").Value = "'6.61%"

# --- Row 8 (MXToken) ---
$ws.Range("D8").Value = "'0.8508"
$ws.Range("E8").Value = "'-0.64%"

# --- Row 9 (FTXToken) ---
$ws.Range("D9").Value = "'0.8633"
$ws.Range("E9").Value = "'-0.69%"

# --- Row 10: was WazirX -> now One ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.0005988"
$ws.Range("E10").Value = "'-0.30%"

# --- Row 11: was MandalaExchangeToken -> now WazirX ---
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1378"
$ws.Range("E11").Value = "'2.27%"

# --- Row 12: was BitrueCoin -> now MandalaExchangeToken ---
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07144"
$ws.Range("E12").Value = "'3.27%"

# --- Row 13: was BitMartToken -> now BitrueCoin ---
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03197"
$ws.Range("E13").Value = "'10.30%"

# --- Row 14: was BitForexToken -> now BitMartToken ---
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09375"
$ws.Range("E14").Value = "'0.00%"

# --- Row 15: was One -> now BitForexToken ---
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001532"
$ws.Range("E15").Value = "'1.63%"

# --- Row 16 (TigerCash) ---
$ws.Range("D16").Value = "'0.005878"
$ws.Range("E16").Value = "'-1.76%"

# --- Row 17 (LEO) ---
$ws.Range("E17").Value = "'-0.26%"

# --- Row 18 (BTSEToken) ---
$ws.Range("D18").Value = "'2.204"
$ws.Range("E18").Value = "'0.94%"

# --- Row 19 (BitpandaEcosystemToken) ---
$ws.Range("E19").Value = "'1.61%"

# --- Row 20 (LiechtensteinCryptoassetsExchange) ---
$ws.Range("D20").Value = "'0.03360"
$ws.Range("E20").Value = "'1.20%"

# --- Row 21 (ProBitToken) ---
$ws.Range("E21").Value = "'-0.33%"

# --- Row 22 (MCDex) ---
$ws.Range("D22").Value = "'3.498"
$ws.Range("E22").Value = "'-2.72%"

# --- Row 23 (CoinExToken) ---
$ws.Range("D23").Value = "'0.04146"
$ws.Range("E23").Value = "'-0.18%"

# --- Row 24 (ZBToken) ---
$ws.Range("D24").Value = "'0.1380"
$ws.Range("E24").Value = "'0.34%"

# --- Row 25 (BitKan) ---
$ws.Range("E25").Value = "'1.32%"

# --- Row 26 (HotbitToken) ---
$ws.Range("E26").Value = "'-7.60%"

# --- Row 27 (NitroEx) ---
$ws.Range("E27").Value = "'1.90%"

# --- Row 28 (UpBots) ---
$ws.Range("D28").Value = "'0.0001448"
$ws.Range("E28").Value = "'4.35%"

# --- Row 40 (IDEX) ---
$ws.Range("E40").Value = "'-0.38%"

# --- Row 41 (KickToken) ---
$ws.Range("D41").Value = "'0.005728"
$ws.Range("E41").Value = "'65.15%"

# --- Row 42 (BKEXToken) ---
$ws.Range("E42").Value = "'0.26%"

# --- Row 43 (CEJI) ---
$ws.Range("E43").Value = "'-3.07%"

# --- Row 44 (LocalTraders) ---
$ws.Range("D44").Value = "'0.009551"
$ws.Range("E44").Value = "'0.06%"

# --- Row 45 (CoinLion) ---
$ws.Range("D45").Value = "'0.00005293"
$ws.Range("E45").Value = "'4.02%"

# --- Row 46 (Kangarootoken) ---
$ws.Range("E46").Value = "'0.20%"

# --- Row 47 (CoinbaseStockToken) ---
$ws.Range("D47").Value = "'0.05798"
$ws.Range("E47").Value = "'-27.36%"

# --- Row 48 (BOLO) ---
$ws.Range("E48").Value = "'-20.46%"

# --- Row 49 (CryptobidCoin) ---
$ws.Range("E49").Value = "'0.20%"

# --- Row 50 (SpecialPowerGold) ---
$ws.Range("E50").Value = "'0.20%"
